# Generate Report for handoff
#
# The localization file c27a0754-b843-451e-a415-02c2846fb04b.md has been
# re-queued for handoff (for both the zh-cn and de-de targets): its status
# flips from "Handed back: in sync with en-US" to "Ready for handoff" and a
# fresh "Latest Handoff Datetime" is recorded per-locale.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"

# --- Overview sheet: summary status column for each locale ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status

# --- zh-cn sheet: detail row for the c27a0754 file ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $status
$zhcn.Range("D3").Value = "2016-01-13 08:31:33"

# --- de-de sheet: detail row for the c27a0754 file ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $status
$dede.Range("D3").Value = "2016-01-13 08:31:53"
